# ---------------------------------------------------------------------------
# da-vinci-ways (Version 2): add a "Meta description" paragraph under the
# H1 title, drop the stray duplicate title paragraph near the end of the
# document, and replace the italic blurb after it with the new AI image
# prompt text.
# ---------------------------------------------------------------------------
$d = $word.ActiveDocument

# --- 1. Insert the "Meta description" paragraph right after the title -----
$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range
$titleRange.Collapse(0)          # wdCollapseEnd
$titleRange.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaRange = $metaPara.Range
$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:r/>' + `
  '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' + `
  '<w:r><w:t>: Experience the beauty of Da Vinci' + [char]39 + 's art while potentially winning big in this slot game with Tumbling Reels and Free Spins. Play for free now.</w:t></w:r>' + `
  '</w:p>'
$metaRange.InsertXML($metaXml) | Out-Null

Write-Output "Meta description paragraph inserted; text now: $($d.Paragraphs(2).Range.Text)"

# --- 2. Remove the stray duplicate "Play Da Vinci Ways..." title paragraph
#        that sits right before the closing italic blurb at the very end --
$count = $d.Paragraphs.Count
$dupPara = $d.Paragraphs($count - 1)
if ($dupPara.Range.Text.TrimEnd() -eq "Play Da Vinci Ways for Free - Slot Game Review") {
    $dupPara.Range.Delete()
}

Write-Output "Count after removal: $($d.Paragraphs.Count)"

# --- 3. Replace the closing italic blurb's wording with the new prompt ----
$newBlurb = 'Please create a cartoon style feature image for the game "Da Vinci Ways". The image should feature a happy Maya warrior with glasses. The warrior should be depicted standing in front of a golden painting that contains the reels of the game, with iconic works of Leonardo da Vinci visible in the painting. The warrior should have a joyful expression on their face, holding up a mobile device with the game on it as if they just won a big payout. The image should convey a sense of fun and excitement while also highlighting the game' + [char]39 + 's theme of art and culture.'

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastFull = $lastPara.Range
$lastText = $d.Range($lastFull.Start, $lastFull.End - 1)   # exclude the paragraph mark
$lastText.Text = $newBlurb

Write-Output "Blurb replaced; last paragraph text now: $($d.Paragraphs($d.Paragraphs.Count).Range.Text)"
